$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1749.5
$ws.Range("I18").Value = 1999.3334
$ws.Range("J18").Value = 1000
$ws.Range("K18").Value = 1999.3334
$ws.Range("L18").Value = 1000
$ws.Range("M18").Value = -1715.3334
$ws.Range("N18").Value = -1568
$ws.Range("H28").Value = 2436.9524
$ws.Range("I28").Value = 1740.5454
$ws.Range("J28").Value = 3203
$ws.Range("K28").Value = 1740.5454
$ws.Range("L28").Value = 3203
$ws.Range("M28").Value = -1255.5454
$ws.Range("N28").Value = -4173
$ws.Range("H48").Value = 4499.5
$ws.Range("J48").Value = 3999
$ws.Range("L48").Value = 11997
$ws.Range("N48").Value = -12581
$ws.Range("H56").Value = 4499.5
$ws.Range("J56").Value = 3999
$ws.Range("L56").Value = 11997
$ws.Range("N56").Value = -13065
$ws.Range("H76").Value = 8486
$ws.Range("I76").Value = 8486
$ws.Range("K76").Value = 8486
$ws.Range("M76").Value = -8171
$ws.Range("H79").Value = 8486
$ws.Range("I79").Value = 8486
$ws.Range("K79").Value = 8486
$ws.Range("M79").Value = -7394
$ws.Range("H116").Value = 4073.4
$ws.Range("I116").Value = 5185
$ws.Range("J116").Value = 3332.3333
$ws.Range("K116").Value = 5185
$ws.Range("L116").Value = 3332.3333
$ws.Range("M116").Value = -1743
$ws.Range("N116").Value = -10216.3333
$ws.Range("H137").Value = 2802.2632
$ws.Range("I137").Value = 2092.5417
$ws.Range("J137").Value = 3318.4243
$ws.Range("K137").Value = 6277.625100000001
$ws.Range("L137").Value = 9955.2729
$ws.Range("M137").Value = -3727.625100000001
$ws.Range("N137").Value = -15055.2729

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1016.3333
$ws.Range("I2").Value = 1016.3333
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1016.3333
$ws.Range("L2").Value = 0
$ws.Range("N2").Value = -903.3333
$ws.Range("H32").Value = 3671.3333
$ws.Range("I32").Value = 1986.3256
$ws.Range("K32").Value = 1986.3256
$ws.Range("M32").Value = -1699.3256
$ws.Range("H43").Value = 68975.14
$ws.Range("I43").Value = 164542.5
$ws.Range("J43").Value = 30748.2
$ws.Range("K43").Value = 164542.5
$ws.Range("L43").Value = 30748.2
$ws.Range("M43").Value = -164229.5
$ws.Range("N43").Value = -31374.2
$ws.Range("H74").Value = 2988.5454
$ws.Range("I74").Value = 2985.875
$ws.Range("J74").Value = 2995.6667
$ws.Range("K74").Value = 2985.875
$ws.Range("L74").Value = 2995.6667
$ws.Range("M74").Value = -2111.875
$ws.Range("N74").Value = -4743.6667
$ws.Range("H77").Value = 2988.5454
$ws.Range("I77").Value = 2985.875
$ws.Range("J77").Value = 2995.6667
$ws.Range("K77").Value = 14929.375
$ws.Range("L77").Value = 14978.3335
$ws.Range("M77").Value = -10561.375
$ws.Range("N77").Value = -23714.3335
$ws.Range("H116").Value = 1016.3333
$ws.Range("I116").Value = 1016.3333
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1016.3333
$ws.Range("L116").Value = 0
$ws.Range("N116").Value = 1277.6667
$ws.Range("H122").Value = 3014.7693
$ws.Range("I122").Value = 3099.4736
$ws.Range("K122").Value = 9298.4208
$ws.Range("M122").Value = -6848.4208

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1016.3333
$ws.Range("I3").Value = 1016.3333
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1016.3333
$ws.Range("L3").Value = 0
$ws.Range("N3").Value = -902.3333
$ws.Range("H99").Value = 1866
$ws.Range("I99").Value = 970.8570999999999
$ws.Range("K99").Value = 970.8570999999999
$ws.Range("M99").Value = 527.1429000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 416.4
$ws.Range("I22").Value = 426.69232
$ws.Range("K22").Value = 426.69232
$ws.Range("M22").Value = -76.69232
$ws.Range("H31").Value = 2647.7188
$ws.Range("J31").Value = 3064
$ws.Range("L31").Value = 3064
$ws.Range("N31").Value = -3654
$ws.Range("H34").Value = 2647.7188
$ws.Range("J34").Value = 3064
$ws.Range("L34").Value = 3064
$ws.Range("N34").Value = -3468
$ws.Range("H122").Value = 4337.9
$ws.Range("I122").Value = 2983
$ws.Range("J122").Value = 4918.5713
$ws.Range("K122").Value = 8949
$ws.Range("L122").Value = 14755.7139
$ws.Range("M122").Value = -6499
$ws.Range("N122").Value = -19655.7139

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 359.5
$ws.Range("J7").Value = 100
$ws.Range("L7").Value = 300
$ws.Range("N7").Value = -524
$ws.Range("H17").Value = 283.33334
$ws.Range("I17").Value = 360
$ws.Range("K17").Value = 1080
$ws.Range("M17").Value = -911
$ws.Range("H113").Value = 675.73334
$ws.Range("I113").Value = 571
$ws.Range("J113").Value = 745.55554
$ws.Range("K113").Value = 1713
$ws.Range("L113").Value = 2236.66662
$ws.Range("M113").Value = 457
$ws.Range("N113").Value = -6576.66662
$ws.Range("H121").Value = 111750.445
$ws.Range("I121").Value = 250340.5
$ws.Range("J121").Value = 878.4
$ws.Range("K121").Value = 751021.5
$ws.Range("L121").Value = 2635.2
$ws.Range("M121").Value = -749711.5
$ws.Range("N121").Value = -5255.2
$ws.Range("H139").Value = 12136.9
$ws.Range("J139").Value = 19962.666
$ws.Range("L139").Value = 59887.99800000001
$ws.Range("N139").Value = -70167.99800000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2336.7693
$ws.Range("I102").Value = 2031.9166
$ws.Range("K102").Value = 2031.9166
$ws.Range("M102").Value = -409.9166
$ws.Range("H122").Value = 2915.8462
$ws.Range("I122").Value = 2900.6365
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 8701.9095
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -6251.9095
$ws.Range("N122").Value = -13898.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 590.6429000000001
$ws.Range("I55").Value = 571.4
$ws.Range("J55").Value = 638.75
$ws.Range("K55").Value = 571.4
$ws.Range("L55").Value = 638.75
$ws.Range("M55").Value = -398.4
$ws.Range("N55").Value = -984.75
$ws.Range("H68").Value = 2663.0715
$ws.Range("I68").Value = 1960.75
$ws.Range("K68").Value = 1960.75
$ws.Range("M68").Value = -1211.75
$ws.Range("H71").Value = 2663.0715
$ws.Range("I71").Value = 1960.75
$ws.Range("K71").Value = 9803.75
$ws.Range("M71").Value = -6059.75
$ws.Range("H108").Value = 75741.664
$ws.Range("J108").Value = 75741.664
$ws.Range("L108").Value = 75741.664
$ws.Range("N108").Value = -83421.664
$ws.Range("H120").Value = 104499.5
$ws.Range("J120").Value = 104499.5
$ws.Range("L120").Value = 104499.5
$ws.Range("N120").Value = -114175.5
$ws.Range("H123").Value = 67500
$ws.Range("J123").Value = 67500
$ws.Range("L123").Value = 67500
$ws.Range("N123").Value = -77300
$ws.Range("H129").Value = 67500
$ws.Range("J129").Value = 67500
$ws.Range("L129").Value = 67500
$ws.Range("N129").Value = -77500
$ws.Range("H130").Value = 63833
$ws.Range("J130").Value = 63833
$ws.Range("L130").Value = 63833
$ws.Range("N130").Value = -73873
$ws.Range("H131").Value = 89500
$ws.Range("J131").Value = 89500
$ws.Range("L131").Value = 89500
$ws.Range("N131").Value = -99580

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 939.26086
$ws.Range("I107").Value = 837.41174
$ws.Range("K107").Value = 2512.23522
$ws.Range("M107").Value = -592.23522
$ws.Range("H108").Value = 69500
$ws.Range("J108").Value = 69500
$ws.Range("L108").Value = 69500
$ws.Range("N108").Value = -77180
$ws.Range("H109").Value = 89900
$ws.Range("J109").Value = 89900
$ws.Range("L109").Value = 89900
$ws.Range("N109").Value = -92674
$ws.Range("H123").Value = 48964.5
$ws.Range("J123").Value = 48964.5
$ws.Range("L123").Value = 48964.5
$ws.Range("N123").Value = -58764.5
